# Auto-generated edit script for "Update latest output (run 39)"
$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates (rows 2-4) ---
$wsSchedule.Range("E2").Value = 1115.13373725
$wsSchedule.Range("F2").Value = 24.58407709986773
$wsSchedule.Range("E3").Value = 362.944725
$wsSchedule.Range("F3").Value = 24.00428075396826
$wsSchedule.Range("A4").Value = 46039.33333333334
$wsSchedule.Range("B4").Value = 46039.83333333334
$wsSchedule.Range("E4").Value = -45.41884425000002
$wsSchedule.Range("F4").Value = -1.001297271825397

# --- Detailed sheet: numeric Price (column B) updates ---
$wsDetailed.Range("B31").Value = 36.0601
$wsDetailed.Range("B32").Value = 28.73596
$wsDetailed.Range("B33").Value = -17.14569
$wsDetailed.Range("B34").Value = -11.91363
$wsDetailed.Range("B35").Value = -7.76371
$wsDetailed.Range("B36").Value = -6.71925
$wsDetailed.Range("B37").Value = -3.18807
$wsDetailed.Range("B38").Value = -0.13614
$wsDetailed.Range("B39").Value = 0.08790000000000001
$wsDetailed.Range("B40").Value = 10.54681
$wsDetailed.Range("B41").Value = 23.33001
$wsDetailed.Range("B42").Value = 32.8827
$wsDetailed.Range("B43").Value = 24.50413
$wsDetailed.Range("B44").Value = 21.974
$wsDetailed.Range("B45").Value = 64.8901
$wsDetailed.Range("B46").Value = 56.98
$wsDetailed.Range("B47").Value = 56.98
$wsDetailed.Range("B54").Value = 36.05988
$wsDetailed.Range("B59").Value = 58.42372
$wsDetailed.Range("B60").Value = 57.06017
$wsDetailed.Range("B62").Value = 56.98
$wsDetailed.Range("B64").Value = 36.06
$wsDetailed.Range("B65").Value = 23.90219
$wsDetailed.Range("B66").Value = -0.8980900000000001
$wsDetailed.Range("B67").Value = -5.33054
$wsDetailed.Range("B68").Value = -2.83936
$wsDetailed.Range("B69").Value = -0.89852
$wsDetailed.Range("B70").Value = -2.54301
$wsDetailed.Range("B71").Value = 0.00886
$wsDetailed.Range("B72").Value = 0.51003
$wsDetailed.Range("B73").Value = 0.51003
$wsDetailed.Range("B74").Value = 0.008630000000000001
$wsDetailed.Range("B75").Value = 6.90404
$wsDetailed.Range("B76").Value = -4.81333
$wsDetailed.Range("B77").Value = -7
$wsDetailed.Range("B78").Value = -11.16992
$wsDetailed.Range("B79").Value = -12.08785
$wsDetailed.Range("B80").Value = -10.065
$wsDetailed.Range("B81").Value = -5.92668
$wsDetailed.Range("B82").Value = 5.27491
$wsDetailed.Range("B83").Value = -9.060980000000001
$wsDetailed.Range("B84").Value = -6.69718
$wsDetailed.Range("B85").Value = -3.05417
$wsDetailed.Range("B86").Value = -2.9124
$wsDetailed.Range("B87").Value = -2.92257
$wsDetailed.Range("B88").Value = 6.81655
$wsDetailed.Range("B90").Value = 57.16514
$wsDetailed.Range("B91").Value = 32.64304
$wsDetailed.Range("B92").Value = 8.499320000000001

# --- Detailed sheet: text (Type / Pump_Status) updates ---
$wsDetailed.Range("C34").Value = "historical"
$wsDetailed.Range("E65").Value = "OFF"
$wsDetailed.Range("E89").Value = "ON"

Write-Output "Done applying updates"
